# Edit script: append the "Five Guys" Team Meeting 5 log section after the
# final paragraph of the document, matching the authored commit
# ("updated meetings log + bug list").

$d = $word.ActiveDocument

# Remember how many paragraphs existed before the insertion so the newly
# created ones can be found afterwards by absolute index.
$baseCount = $d.Paragraphs.Count

# --- 1. Insert all of the new paragraph text in one shot -------------------
# Each "+[char]13+" is a paragraph mark (Word's Enter key). Collapsing the
# last paragraph's range to its end and inserting there places the new text
# immediately before the trailing _GoBack bookmark, exactly like a user
# continuing to type at the end of the document. The leading [char]13 closes
# out the existing last paragraph untouched before the 35 new paragraphs of
# text begin.
$lastPara = $d.Paragraphs.Last
$insertRange = $lastPara.Range
$insertRange.Collapse(0)
$insertRange.InsertAfter([char]13 + "" + [char]13 + "" + [char]13 + "" + [char]13 + "" + [char]13 + "“Five Guys” Team Meeting 5 – 11/18/19" + [char]13 + "" + [char]13 + "ATTENDEES: `tRob, Qui, Afnan, Cameron, Beau, Cole" + [char]13 + "LOCATION: `t`tSpahr Auditorium (Eaton 2)" + [char]13 + "" + [char]13 + "Practice demo went spectacularly, Dr. Gibbons gave us the “smiley face of approval”" + [char]13 + "" + [char]13 + "still to be done" + [char]13 + "progmon switching for both players" + [char]13 + "final boss progmon" + [char]13 + "end screen" + [char]13 + "test suite" + [char]13 + "fourth bag item" + [char]13 + "sound effects" + [char]13 + "hit markers" + [char]13 + "updated documentation" + [char]13 + "product backlog" + [char]13 + "html documentation" + [char]13 + "UML diagrams" + [char]13 + "state diagrams" + [char]13 + "use-case diagrams" + [char]13 + "class diagrams" + [char]13 + "meeting logs" + [char]13 + "gantt chart" + [char]13 + "project 4 tasks" + [char]13 + "bug list" + [char]13 + "users manual" + [char]13 + "deployment plan" + [char]13 + "maintenance plan" + [char]13 + "" + [char]13 + "goal is to have all code finished by SATURDAY of this week so that we can begin work on all of the documentation and not be in a rush to get it completed")

# --- 2. Re-apply correct paragraph formatting ------------------------------
# The freshly typed paragraphs all inherit the formatting of the paragraph
# they were typed from (List Paragraph / list level 2). Walk the newly
# created paragraphs and reset each one to what the log actually needs:
#   plain  -> Normal style, no numbering (blank separator lines, headers)
#   center -> Normal style, no numbering, centered (meeting title block)
#   list   -> List Paragraph style, numId 1, at the given list level
$table = @(
    @(0,'plain',0),
    @(1,'plain',0),
    @(2,'plain',0),
    @(3,'plain',0),
    @(4,'center',0),
    @(5,'center',0),
    @(6,'plain',0),
    @(7,'plain',0),
    @(8,'plain',0),
    @(9,'list',1),
    @(10,'plain',0),
    @(11,'list',1),
    @(12,'list',2),
    @(13,'list',2),
    @(14,'list',2),
    @(15,'list',2),
    @(16,'list',2),
    @(17,'list',2),
    @(18,'list',2),
    @(19,'list',2),
    @(20,'list',3),
    @(21,'list',3),
    @(22,'list',3),
    @(23,'list',4),
    @(24,'list',4),
    @(25,'list',4),
    @(26,'list',3),
    @(27,'list',3),
    @(28,'list',4),
    @(29,'list',3),
    @(30,'list',3),
    @(31,'list',3),
    @(32,'list',3),
    @(33,'plain',0),
    @(34,'list',1)
)

foreach ($row in $table) {
    $offset = $row[0]
    $kind = $row[1]
    $level = $row[2]
    $para = $d.Paragraphs.Item($baseCount + $offset + 1)
    $rng = $para.Range

    if ($kind -eq 'list') {
        $rng.ListFormat.ListLevelNumber = $level
    } else {
        $rng.ListFormat.RemoveNumbers()
        $para.Style = $d.Styles.Item("Normal")
        if ($kind -eq 'center') {
            $para.Alignment = 1
        } else {
            $para.Alignment = 0
        }
    }
}

Write-Output "inserted $($table.Count) paragraphs after index $baseCount"
